$d = $word.ActiveDocument

$find = "Dates de la campanya 2022 en què usem la constel·lació, Constel·lació de Taure 16-25 de gener"
$replace = "Dates de la campanya 2022 en què usem la  Constel·lació de Taure 16-25 de gener"

$rng = $d.Content
$rng.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
